$wb = $excel.ActiveWorkbook

# Week 16 logged on the OFF sheet (row 2 = "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 352
$wsOff.Range("C2").Value = 228
$wsOff.Range("D2").Value = 149
$wsOff.Range("E2").Value = 65

# Week 16 logged on the DEF sheet (row 2 = "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 483
$wsDef.Range("C2").Value = 364
$wsDef.Range("D2").Value = 95
$wsDef.Range("E2").Value = 57
